$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E3").Value = 34
$ws.Range("E7").Value = 28
$ws.Range("E10").Value = 414
$ws.Range("F10").Value = 207
$ws.Range("H10").Value = 207
$ws.Range("E11").Value = 279
$ws.Range("F11").Value = 153
$ws.Range("H11").Value = 153
$ws.Range("E12").Value = 405
$ws.Range("F12").Value = 221
$ws.Range("H12").Value = 221
$ws.Range("E13").Value = 107
$ws.Range("E15").Value = 134
$ws.Range("E16").Value = 172
$ws.Range("F16").Value = 89
$ws.Range("H16").Value = 89
$ws.Range("E20").Value = 79
$ws.Range("E21").Value = 125
$ws.Range("E23").Value = 174
$ws.Range("F23").Value = 76
$ws.Range("H23").Value = 76
$ws.Range("E24").Value = 177
$ws.Range("F24").Value = 95
$ws.Range("H24").Value = 95
$ws.Range("E25").Value = 214
$ws.Range("E26").Value = 121
$ws.Range("E27").Value = 278
$ws.Range("F27").Value = 131
$ws.Range("H27").Value = 131
$ws.Range("E28").Value = 165
$ws.Range("E33").Value = 249
$ws.Range("E34").Value = 180
$ws.Range("E35").Value = 118
$ws.Range("F35").Value = 77
$ws.Range("H35").Value = 77
$ws.Range("E36").Value = 55
$ws.Range("F36").Value = 33
$ws.Range("H36").Value = 33
$ws.Range("E40").Value = 225
$ws.Range("F40").Value = 103
$ws.Range("H40").Value = 103
$ws.Range("E41").Value = 333
$ws.Range("F41").Value = 156
$ws.Range("H41").Value = 156
$ws.Range("E42").Value = 303
$ws.Range("E43").Value = 101
$ws.Range("E44").Value = 259
$ws.Range("F44").Value = 127
$ws.Range("H44").Value = 127
$ws.Range("E45").Value = 117
$ws.Range("F45").Value = 58
$ws.Range("H45").Value = 58
$ws.Range("F46").Value = 143
$ws.Range("H46").Value = 143
$ws.Range("E47").Value = 375
$ws.Range("E48").Value = 176
$ws.Range("F48").Value = 71
$ws.Range("H48").Value = 71
$ws.Range("E49").Value = 252
$ws.Range("E50").Value = 217
$ws.Range("F50").Value = 94
$ws.Range("H50").Value = 94
$ws.Range("E51").Value = 201
$ws.Range("F51").Value = 80
$ws.Range("H51").Value = 80
